$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SARIMAX")
$ws1.Activate()
$ws1.Range("A88").Select()
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("F113").Select()
